$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Skeena")

# Add a new fisheries-opening entry for Aug 23-29 (FN0853) on row 37.
$ws.Range("A37").Value = "FN0853"
$ws.Range("B37").Value = "Aboriginal"
$ws.Range("C37").Value = "Aug 23-29"
$ws.Range("C37").NumberFormat = "d-mmm"
$ws.Range("D37").Value = "Sockeye"
$ws.Range("E37").Value = "Selective Gear"
$ws.Range("F37").Value = "Region 6-Lake Babine Nation"
$ws.Range("G37").Value = 7
$ws.Range("I37").Value = "Sockeye retention only"

# Update the "Region 6-LBN" label on row 35 (Aug 15-21 entry) to its fuller name.
$ws.Range("F35").Value = "Region 6-Lake Babine Nation"
